$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 52, shifting the existing rows 52:101 down to 53:102.
$ws.Rows.Item(52).Insert()

# Populate the newly inserted row 52 with a new weekly price record.
# Static columns are identical across all rows in this data block.
$ws.Cells.Item(52, 1).Value = 5
$ws.Cells.Item(52, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(52, 3).Value = "Maule"
$ws.Cells.Item(52, 4).Value = 44589
$ws.Cells.Item(52, 5).Value = 7
$ws.Cells.Item(52, 6).Value = 100112030
$ws.Cells.Item(52, 7).Value = "Poroto granado"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 300
$ws.Cells.Item(52, 11).Value = 25000
$ws.Cells.Item(52, 12).Value = 25000
$ws.Cells.Item(52, 13).Value = 25000
$ws.Cells.Item(52, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(52, 15).Value = "Región del Maule"
$ws.Cells.Item(52, 16).Value = 1000
$ws.Cells.Item(52, 17).Value = 25
$ws.Cells.Item(52, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Range("D52").NumberFormat = $ws.Range("D53").NumberFormat
